# Automatische test-sync: 2025-07-23 21:07:50
# Appends a new testmail log row to "Logs" and bumps the matching
# "Dashboard" counter, extending conditional formatting ranges to cover
# the newly added row.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append row 8 to the Logs sheet -----------------------------------
$logs.Range("A8").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Range("D8").Value = "Openingstijden / Locatie"
$logs.Range("E8").Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F8").Value = "2025-07-23 21:07:23"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"
$logs.Range("J8").Value = "Ja"

# Writing the multi-line E8 text auto-pins a custom row height; re-run
# AutoFit so the row reverts to the default (no ht/customHeight attrs),
# matching the source rows which never had an explicit height either.
$logs.Rows.Item(8).EntireRow.AutoFit()

# --- Extend conditional formatting ranges to include row 8 -----------
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))

# --- Bump the Dashboard tally for this category -----------------------
$dashboard.Range("B2").Value = 7
